$wb = $excel.ActiveWorkbook

# The "最低票价" (minimum ticket price) column G was re-scaled from
# fen-like integers (x100) down to yuan (/100) and is now stored as text
# rather than a number. Writing a leading apostrophe forces the
# numeric-looking text to be kept as a string instead of being re-parsed
# as a number; ClearFormats() then drops the transient "quote prefix"
# cell style Excel applies for that so the cell keeps its original
# (default) style, just like the rest of the sheet.
function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: G2 5400 -> "54" (text)
    Set-TextValue $ws.Range("G2") "54"

    # Row 3: G3 5000 -> "50" (text)
    Set-TextValue $ws.Range("G3") "50"

    # Row 4: F4 1430 -> 1431 ; G4 6000 -> "60" (text)
    $ws.Range("F4").Value = 1431
    Set-TextValue $ws.Range("G4") "60"

    # Row 5: F5 6948 -> 6951 ; G5 5500 -> "55" (text)
    $ws.Range("F5").Value = 6951
    Set-TextValue $ws.Range("G5") "55"

    # Row 6: G6 5000 -> "50" (text)
    Set-TextValue $ws.Range("G6") "50"

    # Row 7: F7 103 -> 104 ; G7 6000 -> "60" (text)
    $ws.Range("F7").Value = 104
    Set-TextValue $ws.Range("G7") "60"
}
